{"js": "// Remove the \"Appendix: Quick prototype\" sub-section (heading, the blank\n// paragraph under it, the \"Figure: PDF page 1\" caption, and the paragraph\n// holding the embedded prototype screenshot) that sits between the\n// \"Appendix: Links\" intro paragraph and the real \"Appendix: Links\" heading.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\n// Locate the \"Appendix: Quick prototype\" Heading 2 paragraph. (There is\n// also an unrelated Normal-style paragraph earlier in the document that\n// starts with the same words, so match on exact paragraph text + style to\n// avoid grabbing it.)\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text === \"Appendix: Quick prototype\" && para.style === \"Heading 2\") {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex !== -1) {\n  // The section spans 4 paragraphs: the Heading 2, a blank paragraph, the\n  // \"Figure: PDF page 1\" caption, and the paragraph with the picture.\n  // Delete from the last one back to the first so earlier indices stay valid.\n  for (let i = targetIndex + 3; i >= targetIndex; i--) {\n    paragraphs.items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the \"Appendix: Quick prototype\" sub-section (heading, the blank\n# paragraph under it, the \"Figure: PDF page 1\" caption, and the paragraph\n# holding the embedded prototype screenshot) that sits between the\n# \"Appendix: Links\" intro paragraph and the real \"Appendix: Links\" heading.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Appendix: Quick prototype\" Heading 2 paragraph (there is also\n# a Normal-style paragraph earlier in the doc that starts with the same\n# words, so match on exact paragraph text + style to avoid grabbing it).\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.TrimEnd(\"`r\")\n    if ($txt -eq \"Appendix: Quick prototype\" -and $p.Style.NameLocal -eq \"Heading 2\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -ne -1) {\n    # The section spans 4 paragraphs: the Heading 2, a blank paragraph,\n    # the \"Figure: PDF page 1\" caption, and the paragraph with the picture.\n    $startPara = $d.Paragraphs.Item($targetIndex)\n    $endPara = $d.Paragraphs.Item($targetIndex + 3)\n\n    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $delRange.Delete()\n}\n"}
